# "configurazione br e moduli nuovi 18-5-2017"
# Add a new analysis-unit-variable row for RETAIL_IND_119 /
# COUNTERPARTY_RETAIL_IND_119 on the "r AnalysisUnit_Variable" sheet,
# just above the existing IND_150 row (new row 99; old rows 99-101 shift
# down to 100-102).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)
$ws.Activate()

# Push the existing rows 99-101 down to 100-102 and open up a blank row 99.
$ws.Rows.Item(99).Insert()

# Fill the new row. Column F ("RETAIL_IND_119") is written before columns
# B/C ("COUNTERPARTY_RETAIL_IND_119") so the shared-string table grows in
# the same order as the source edit (RETAIL_IND_119 then
# COUNTERPARTY_RETAIL_IND_119).
$ws.Range("A99").Value = "CREATE/MODIFY"
$ws.Range("F99").Value = "RETAIL_IND_119"
$ws.Range("B99").Value = "COUNTERPARTY_RETAIL_IND_119"
$ws.Range("C99").Value = "COUNTERPARTY_RETAIL_IND_119"
$ws.Range("E99").Value = "COUNTERPARTY_RETAIL"

# Match the author's final view/selection state on the sheet.
$ws.Range("E99").Select()
